# Fixed Diagnosis, FileAssociation, FileFormat, FileType, NeuteredStatus, PrimeDiseaseSite
#
# The "CasesTab" Cypher query (stored in cell B2 of the "startup" sheet) is
# updated to drop the trailing `Cohort` column from its RETURN clause, since
# the cohort node is no longer projected by that query.

$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("startup")

$newCasesQuery = @'
MATCH (s:study)<-[*]-(c:case)<--(demo:demographic)
MATCH (c)<--(diag:diagnosis)
OPTIONAL MATCH (samp:sample)-->(c)
OPTIONAL MATCH (co:cohort)<-[*]-(c)
WITH DISTINCT c, s, demo, diag, co
WHERE demo.neutered_indicator IN ["Yes"]  
RETURN  coalesce(c.case_id, '') AS `Case ID` ,
        coalesce(s.clinical_study_designation, '') AS `Study Code` ,
        coalesce(s.clinical_study_type, '') AS  `Study Type`,
        coalesce(demo.breed, '') AS Breed ,
        coalesce(diag.disease_term, '') AS Diagnosis ,
        coalesce(diag.stage_of_disease, '') AS `Stage of Disease` ,
        coalesce(demo.patient_age_at_enrollment, '') AS Age ,
        coalesce(demo.sex, '') AS Sex ,
        coalesce(demo.neutered_indicator, '') AS `Neutered Status`,
        coalesce(demo.weight, '') AS `Weight (kg)`,
        coalesce(diag.best_response, '') AS `Response to Treatment`
'@

$ws.Range("B2").Value = $newCasesQuery

# Restore the selection to B2 (it had moved to B4 previously).
[void]$ws.Range("B2").Select()
